$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: restyle cells that used the removed style (s=3) and still have values ---
# Setting VerticalAlignment to top matches style index 1, reusing it and freeing style index 3.
$restyleCells = @("F199","F200","F201","F202","I202","F203","I203","F206","F207","I207","F210","I210","I212")
foreach ($addr in $restyleCells) {
  $ws.Range($addr).VerticalAlignment = -4160
}

# --- Step 2: fully remove cells that used the removed style (s=3) but carried no value ---
$clearCells = @("I204","I208","F209","I209","F211","I211")
foreach ($addr in $clearCells) {
  $ws.Range($addr).Clear()
}

# --- Step 3: append new sprint S24 rows (213-225) ---
$ws.Range("A213").Value = 'S24'
$ws.Range("B213").Value = 'G01'
$ws.Range("C213").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D213").Value = 'S24_G01_TB001'
$ws.Range("E213").Value = 'Backend: Add DB tables/models for group import datasets (schema + per-symbol values) and link dataset to a group.'
$ws.Range("G213").Value = 'implemented'
$ws.Range("H213").Value = 'Added GroupImport/GroupImportValue models + Alembic migration 0032.'
$ws.Rows.Item(213).RowHeight = 41.75

$ws.Range("A214").Value = 'S24'
$ws.Range("B214").Value = 'G01'
$ws.Range("C214").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D214").Value = 'S24_G01_TB002'
$ws.Range("E214").Value = 'Backend: Implement symbol normalization + broker instrument resolution for NSE/BSE; skip/reject rows that do not resolve.'
$ws.Range("G214").Value = 'implemented'
$ws.Range("H214").Value = 'Implemented symbol normalization + broker instrument validation using market_instruments (+ optional Kite fallback).'
$ws.Rows.Item(214).RowHeight = 41.75

$ws.Range("A215").Value = 'S24'
$ws.Range("B215").Value = 'G01'
$ws.Range("C215").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D215").Value = 'S24_G01_TB003'
$ws.Range("E215").Value = 'Backend: Enforce import column restrictions (no OHLCV/price/volume/performance/indicator-derived columns); return polite reasons for skipped columns.'
$ws.Range("G215").Value = 'implemented'
$ws.Range("H215").Value = 'Enforced disallowed import columns (OHLCV/perf/indicators/ratios) in backend + UI.'
$ws.Rows.Item(215).RowHeight = 41.75

$ws.Range("A216").Value = 'S24'
$ws.Range("B216").Value = 'G01'
$ws.Range("C216").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D216").Value = 'S24_G01_TB004'
$ws.Range("E216").Value = 'Backend: Add group import API endpoint to create group + dataset from uploaded/parsed CSV with mapping/selection instructions; support “replace vs keep existing dataset”.'
$ws.Range("G216").Value = 'implemented'
$ws.Range("H216").Value = 'Added POST /api/groups/import/watchlist with conflict handling (ERROR/REPLACE_DATASET/REPLACE_GROUP).'
$ws.Rows.Item(216).RowHeight = 55.2

$ws.Range("A217").Value = 'S24'
$ws.Range("B217").Value = 'G01'
$ws.Range("C217").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D217").Value = 'S24_G01_TB005'
$ws.Range("E217").Value = 'Backend: Add read endpoint(s) to return dataset schema + values for a group so Holdings/Groups grid can render dynamic columns.'
$ws.Range("G217").Value = 'implemented'
$ws.Range("H217").Value = 'Added GET /api/groups/{id}/dataset and /dataset/values for dynamic column rendering.'
$ws.Rows.Item(217).RowHeight = 41.75

$ws.Range("A218").Value = 'S24'
$ws.Range("B218").Value = 'G01'
$ws.Range("C218").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D218").Value = 'S24_G01_TB006'
$ws.Range("E218").Value = 'Backend: Tests for symbol resolution, restricted columns, and import endpoint behaviors (duplicate headers, collisions, replace semantics).'
$ws.Range("G218").Value = 'implemented'
$ws.Range("H218").Value = 'Added backend tests covering restrictions, unresolved symbols, and replace semantics.'
$ws.Rows.Item(218).RowHeight = 41.75

$ws.Range("A219").Value = 'S24'
$ws.Range("B219").Value = 'G01'
$ws.Range("C219").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D219").Value = 'S24_G01_TF001'
$ws.Range("E219").Value = 'Frontend: Add “Import watchlist (CSV)” wizard on Groups page: upload → preview → map symbol/exchange → choose columns → create group.'
$ws.Range("G219").Value = 'implemented'
$ws.Range("H219").Value = 'Added Groups page CSV import wizard (upload → map → select columns → import).'
$ws.Rows.Item(219).RowHeight = 41.75

$ws.Range("A220").Value = 'S24'
$ws.Range("B220").Value = 'G01'
$ws.Range("C220").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D220").Value = 'S24_G01_TF002'
$ws.Range("E220").Value = 'Frontend: Render dynamic dataset columns in Holdings/Groups grid views only when selected group has an attached dataset.'
$ws.Range("G220").Value = 'implemented'
$ws.Range("H220").Value = 'Holdings group view now loads dataset and renders dynamic DataGrid columns for imported metadata.'
$ws.Rows.Item(220).RowHeight = 41.75

$ws.Range("A221").Value = 'S24'
$ws.Range("B221").Value = 'G01'
$ws.Range("C221").Value = 'Groups improvement: watchlist CSV import with dynamic columns (see docs/groups_improvement.md)'
$ws.Range("D221").Value = 'S24_G01_TF003'
$ws.Range("E221").Value = 'Frontend: Show import summary with skipped symbols/columns + reasons; provide “replace dataset” confirmation for existing group.'
$ws.Range("G221").Value = 'implemented'
$ws.Range("H221").Value = 'Import dialog shows summary (imported/skipped) and links to open imported group in holdings grid.'
$ws.Rows.Item(221).RowHeight = 41.75

$ws.Range("A222").Value = 'S24'
$ws.Range("B222").Value = 'G02'
$ws.Range("C222").Value = 'Groups improvement (phase 2): XLSX support + portfolio import mappings'
$ws.Range("D222").Value = 'S24_G02_TB001'
$ws.Range("E222").Value = 'Backend: Add XLSX import support (parse on backend) and reuse the same dataset creation pipeline as CSV.'
$ws.Range("G222").Value = 'planned'
$ws.Range("I222").Value = 'Ship after CSV flow is stable.'
$ws.Rows.Item(222).RowHeight = 28.35

$ws.Range("A223").Value = 'S24'
$ws.Range("B223").Value = 'G02'
$ws.Range("C223").Value = 'Groups improvement (phase 2): XLSX support + portfolio import mappings'
$ws.Range("D223").Value = 'S24_G02_TF001'
$ws.Range("E223").Value = 'Frontend: Accept .xlsx in import wizard and upload to backend for parsing; keep mapping/selection UX unchanged.'
$ws.Range("G223").Value = 'planned'
$ws.Range("I223").Value = 'Avoid heavy client-side XLSX parsing initially.'
$ws.Rows.Item(223).RowHeight = 41.75

$ws.Range("A224").Value = 'S24'
$ws.Range("B224").Value = 'G02'
$ws.Range("C224").Value = 'Groups improvement (phase 2): XLSX support + portfolio import mappings'
$ws.Range("D224").Value = 'S24_G02_TB002'
$ws.Range("E224").Value = 'Backend: Add portfolio group kind + portfolio fields storage (qty, avg_buy_price, buy_date, target_weight) keyed by (group, instrument).'
$ws.Range("G224").Value = 'planned'
$ws.Range("I224").Value = 'Keep watchlist datasets separate from portfolio fields.'
$ws.Rows.Item(224).RowHeight = 41.75

$ws.Range("A225").Value = 'S24'
$ws.Range("B225").Value = 'G02'
$ws.Range("C225").Value = 'Groups improvement (phase 2): XLSX support + portfolio import mappings'
$ws.Range("D225").Value = 'S24_G02_TF002'
$ws.Range("E225").Value = 'Frontend: Portfolio import mapping UI (map file columns to qty/avg_buy/buy_date/weight) + validation; still allow extra dynamic columns.'
$ws.Range("G225").Value = 'planned'
$ws.Range("I225").Value = 'Phase 2/3 work; requires clear portfolio UX.'
$ws.Rows.Item(225).RowHeight = 41.75

# --- Step 4: update sheet view (scroll position + selection) to match post-edit state ---
$ws.Application.Goto($ws.Range("A215"), $true)
$ws.Cells.Select()
$ws.Range("C216").Select()
